$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("<india>",   "<india>",   8),
    @("<uniform>", "<uniform>", 7),
    @("<oscar>",   "<oscar>",   6),
    @("<water>",   "<water>",   6),
    @("<so>",      "<so>",      5),
    @("<and>",     "<and>",     5),
    @("<zero>",    "<zero>",    7),
    @("<when>",    "<when>",    4),
    @("<sentence>","<sentence>",5),
    @("<cut>",     "<cut>",     4),
    @("<could>",   "<could>",   1),
    @("<delta>",   "<quef>",    13),
    @("<that>",    "<thebec>",  6),
    @("<first>",   "<first>",   5),
    @("<could>",   "<could>",   5),
    @("<then>",    "<then>",    5),
    @("<can>",     "<can>",     6)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row++
}
